# fix: prevent hidden columns from being labeled upon detecting changes
#
# Rows whose only "difference" between the FV2410 and FV2504 columns came
# from a (now-ignored/hidden) column no longer qualify as changed rows, so
# their "Änderung" (= "ÄNDERUNG") marker in column L must be cleared.
#
# A subset of those rows are also the first ("parent"/group header) row of
# a segment-group block; those rows additionally need their whole-row
# highlighting style corrected (same pattern already used correctly on
# e.g. row 9) - grey fill for every cell (style index 6), with the bold
# "Segmentname" cell in column B using style index 7, and the now-empty
# column L cell using style index 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that are the first row of a segment group and need their full
# row style corrected (A:V -> style 6, except B -> style 7, L -> style 5).
$parentRows = @(13, 17, 23, 27, 34, 40, 63, 67)

# Row 9 already carries the exact target style pattern for a parent row,
# so reuse it as the formatting template instead of hand-rolling style
# indices.
$styleTemplateRow = $ws.Range("A9:V9")
$styleTemplateRow.Copy()
foreach ($r in $parentRows) {
    $target = $ws.Range("A$r`:V$r")
    $target.PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = $false

# All rows (parent rows included) whose "ÄNDERUNG" label in column L must
# be removed because the row no longer counts as changed.
$rowsToClearChangeMarker = @(
    13, 14, 15, 16, 17, 18, 19, 20, 21, 22,
    23, 24, 25, 26, 27, 28, 29, 30, 31, 32,
    33, 34, 35, 36, 38, 39, 40, 41, 42, 43,
    63, 67, 104, 115
)

# L9 is already an empty cell using the correct "no marker" style (5), so
# copy its formatting onto every column-L cell that needs to lose its
# marker, then clear the cell's value/content.
$styleTemplateL = $ws.Range("L9")
$styleTemplateL.Copy()
foreach ($r in $rowsToClearChangeMarker) {
    $cell = $ws.Range("L$r")
    $cell.PasteSpecial(-4122)  # xlPasteFormats
    $cell.Value = ""
}
$excel.CutCopyMode = $false
